$d = $word.ActiveDocument

# 1) Date in the daily header line: 07.07.24 -> 05.07.24
$d.Paragraphs(1).Range.Text = "⚡️🚀המאמר היומי של מייק 05.07.24:⚡️🚀"

# 2) Paper title
$d.Paragraphs(2).Range.Text = "A Survey of Large Language Models for Graphs"

# 3) First body paragraph (intro)
$d.Paragraphs(3).Range.Text = "גרפים מודלי שפה גדולים: האם זה שידוך מהחלומות? גרפים נמצאים בכל מקום, מרשתות חברתיות ועד למבנים מולקולריים ורשתות נוירונים על גרפים (GNNs) הם הפתרון הנפוץ למשימות כמו ניבוי קישורים וסיווג קודקודים. אבל ל-GNNs יש מגבלות: הם מתקשים עם דאטה דליל ולעיתים קרובות אינם מצליחים להכליל היטב לגרפים בעל מבנה שלא נראו קודם."

# 4) Second body paragraph
$d.Paragraphs(4).Range.Text = "מאידך גיסא LLMs מספקים פתרון משלים: הם מצטיינים בהבנה וסיכום טקסטים (שזה דאטה דליל שהוא בעצם גרף - המתאר קשרים בין מילים או קבוצות של מילים) יותר מאשר גרפים. אז, מה אם נשלב את החוזקות של GNNs ו-LLMs? מאמר סקר חדש חוקר לעומק את החיבור המבטיח הזה."

# 5) Third body paragraph
$d.Paragraphs(5).Range.Text = "המחברים מציעים טקסונומיה של ארבעה שילובים אפשריים בין LLM ל-GNNs: שימוש ב-GNNs בתור שלב מקדים ל-LLMs, שימוש ב-LLMs לפני GNNs, שילוב של LLMs וגרפים, ושימוש ב-LLMs בלבד למשימות גרפיות. לכל גישה יש יתרונות וחסרונות, אבל הפוטנציאל ברור. על ידי ניצול הכוח של LLMs, נוכל להתגבר על חלק מהמגבלות של טכניקות למידה מסורתיות על גרפים."

# 6) Link
$d.Paragraphs(6).Range.Text = "https://arxiv.org/pdf/2405.08011"
